$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-looking text (column H, HUMITAT_MITJANA_DIA) needs a leading
# quote-prefix (like typing an apostrophe before 55% by hand in Excel) so
# COM stores it as the literal text "55%" instead of auto-converting it to
# the number 0.55.

$ws.Range('E2').Value = '2026-02-07 15:17:38'
$ws.Range('K2').Value = '8.3 MJ/m2'
$ws.Range('O2').Value = '-0.8 °C'
$ws.Range('E3').Value = '2026-02-07 15:17:40'
$ws.Range('K3').Value = '12.8 MJ/m2'
$ws.Range('L3').Value = '23.4 km/h - 104º 14:39 TU'
$ws.Range('E4').Value = '2026-02-07 15:17:42'
$ws.Range('H4').Value = "'55%"
$ws.Range('K4').Value = '9.3 MJ/m2'
$ws.Range('M4').Value = '15.3 °C 14:39 TU'
$ws.Range('O4').Value = '11.9 °C'
$ws.Range('E5').Value = '2026-02-07 15:17:45'
$ws.Range('H5').Value = "'66%"
$ws.Range('J5').Value = '1003.3 hPa'
$ws.Range('K5').Value = '11.0 MJ/m2'
$ws.Range('O5').Value = '10.3 °C'
$ws.Range('E6').Value = '2026-02-07 15:17:47'
$ws.Range('H6').Value = "'48%"
$ws.Range('J6').Value = '1004.9 hPa'
$ws.Range('K6').Value = '11.2 MJ/m2'
$ws.Range('M6').Value = '16.5 °C 14:50 TU'
$ws.Range('O6').Value = '12.9 °C'
$ws.Range('E7').Value = '2026-02-07 15:17:50'
$ws.Range('J7').Value = '1004.4 hPa'
$ws.Range('K7').Value = '11.8 MJ/m2'
$ws.Range('M7').Value = '14.1 °C 14:46 TU'
$ws.Range('O7').Value = '9.1 °C'
$ws.Range('E8').Value = '2026-02-07 15:17:52'
$ws.Range('H8').Value = "'77%"
$ws.Range('K8').Value = '10.6 MJ/m2'
$ws.Range('O8').Value = '8.3 °C'
$ws.Range('E9').Value = '2026-02-07 15:17:55'
$ws.Range('H9').Value = "'88%"
$ws.Range('M9').Value = '11.0 °C 14:37 TU'
$ws.Range('O9').Value = '3.3 °C'
$ws.Range('E10').Value = '2026-02-07 15:17:57'
$ws.Range('O10').Value = '10.2 °C'
$ws.Range('E11').Value = '2026-02-07 15:17:59'
$ws.Range('H11').Value = "'86%"
$ws.Range('J11').Value = '1006.1 hPa'
$ws.Range('K11').Value = '7.7 MJ/m2'
$ws.Range('O11').Value = '3.0 °C'
$ws.Range('E12').Value = '2026-02-07 15:18:01'
$ws.Range('K12').Value = '11.4 MJ/m2'
$ws.Range('O12').Value = '12.1 °C'
$ws.Range('E13').Value = '2026-02-07 15:18:04'
$ws.Range('O13').Value = '11.2 °C'
$ws.Range('E14').Value = '2026-02-07 15:18:06'
$ws.Range('K14').Value = '8.1 MJ/m2'
$ws.Range('E15').Value = '2026-02-07 15:18:09'
$ws.Range('H15').Value = "'73%"
$ws.Range('K15').Value = '10.5 MJ/m2'
$ws.Range('O15').Value = '9.3 °C'
$ws.Range('E16').Value = '2026-02-07 15:18:11'
$ws.Range('H16').Value = "'88%"
$ws.Range('K16').Value = '5.6 MJ/m2'
$ws.Range('M16').Value = '8.9 °C 14:32 TU'
$ws.Range('O16').Value = '3.7 °C'
$ws.Range('E17').Value = '2026-02-07 15:18:14'
$ws.Range('H17').Value = "'89%"
$ws.Range('J17').Value = '1005.8 hPa'
$ws.Range('K17').Value = '7.8 MJ/m2'
$ws.Range('M17').Value = '9.7 °C 14:39 TU'
$ws.Range('O17').Value = '4.6 °C'
$ws.Range('E18').Value = '2026-02-07 15:18:16'
$ws.Range('K18').Value = '5.5 MJ/m2'
$ws.Range('O18').Value = '-5.8 °C'
$ws.Range('E19').Value = '2026-02-07 15:18:18'
$ws.Range('H19').Value = "'84%"
$ws.Range('K19').Value = '11.2 MJ/m2'
$ws.Range('O19').Value = '6.7 °C'
$ws.Range('E20').Value = '2026-02-07 15:18:21'
$ws.Range('K20').Value = '9.7 MJ/m2'
$ws.Range('E21').Value = '2026-02-07 15:18:23'
$ws.Range('J21').Value = '1003.6 hPa'
$ws.Range('K21').Value = '10.5 MJ/m2'
$ws.Range('M21').Value = '13.9 °C 14:46 TU'
$ws.Range('O21').Value = '8.1 °C'
$ws.Range('E22').Value = '2026-02-07 15:18:26'
$ws.Range('H22').Value = "'75%"
$ws.Range('K22').Value = '11.9 MJ/m2'
$ws.Range('O22').Value = '10.1 °C'
$ws.Range('E23').Value = '2026-02-07 15:18:28'
$ws.Range('H23').Value = "'77%"
$ws.Range('J23').Value = '1003.3 hPa'
$ws.Range('K23').Value = '9.4 MJ/m2'
$ws.Range('O23').Value = '10.2 °C'
$ws.Range('E24').Value = '2026-02-07 15:18:30'
$ws.Range('J24').Value = '1002.7 hPa'
$ws.Range('K24').Value = '9.0 MJ/m2'
$ws.Range('O24').Value = '11.0 °C'
$ws.Range('E25').Value = '2026-02-07 15:18:33'
$ws.Range('H25').Value = "'87%"
$ws.Range('I25').Value = '2.1 mm'
$ws.Range('K25').Value = '6.2 MJ/m2'
$ws.Range('M25').Value = '6.9 °C 14:41 TU'
$ws.Range('O25').Value = '1.9 °C'
$ws.Range('E26').Value = '2026-02-07 15:18:35'
$ws.Range('H26').Value = "'68%"
$ws.Range('K26').Value = '10.3 MJ/m2'
$ws.Range('O26').Value = '-1.8 °C'
$ws.Range('E27').Value = '2026-02-07 15:18:38'
$ws.Range('J27').Value = '1003.4 hPa'
$ws.Range('K27').Value = '10.1 MJ/m2'
$ws.Range('O27').Value = '11.3 °C'
$ws.Range('E28').Value = '2026-02-07 15:18:40'
$ws.Range('H28').Value = "'83%"
$ws.Range('J28').Value = '1005.4 hPa'
$ws.Range('O28').Value = '4.3 °C'
$ws.Range('E29').Value = '2026-02-07 15:18:43'
$ws.Range('K29').Value = '11.4 MJ/m2'
$ws.Range('O29').Value = '12.2 °C'
$ws.Range('E30').Value = '2026-02-07 15:18:45'
$ws.Range('K30').Value = '13.0 MJ/m2'
$ws.Range('E31').Value = '2026-02-07 15:18:47'
$ws.Range('H31').Value = "'87%"
$ws.Range('J31').Value = '1006.6 hPa'
$ws.Range('M31').Value = '10.1 °C 14:34 TU'
$ws.Range('O31').Value = '5.3 °C'
$ws.Range('E32').Value = '2026-02-07 15:18:50'
$ws.Range('H32').Value = "'49%"
$ws.Range('K32').Value = '11.4 MJ/m2'
$ws.Range('O32').Value = '12.9 °C'
$ws.Range('E33').Value = '2026-02-07 15:18:52'
$ws.Range('H33').Value = "'78%"
$ws.Range('O33').Value = '10.1 °C'
$ws.Range('E34').Value = '2026-02-07 15:18:54'
$ws.Range('H34').Value = "'73%"
$ws.Range('K34').Value = '9.1 MJ/m2'
$ws.Range('O34').Value = '7.5 °C'
$ws.Range('E35').Value = '2026-02-07 15:18:57'
$ws.Range('K35').Value = '5.4 MJ/m2'
$ws.Range('O35').Value = '-4.5 °C'
$ws.Range('E36').Value = '2026-02-07 15:18:59'
$ws.Range('H36').Value = "'80%"
$ws.Range('K36').Value = '11.0 MJ/m2'
$ws.Range('O36').Value = '8.1 °C'
